# Updates the "Price" (column D) and "Volume(1h)" (column E) values for the
# cryptocurrency rows on Sheet1, matching the latest scrape from the
# GitHub Actions job that refreshes cryptos.xlsx.
#
# Each row's Price/Volume text is only rewritten if the diff actually
# changed it for that row (some rows only changed one of the two columns,
# or neither).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "29.348.45";    E = "  +0.15%  " },
    @{ Row = 3;  D = "1.878.07";     E = "  +0.24%  " },
    @{ Row = 4;  D = "1.0000";       E = "  +0.01%  " },
    @{ Row = 5;  D = "0.7146";       E = "  +0.35%  " },
    @{ Row = 6;  D = "242.24";       E = "  -0.18%  " },
    @{ Row = 7;  D = $null;          E = "  +0.08%  " },
    @{ Row = 8;  D = "0.08113";      E = "  +4.69%  " },
    @{ Row = 9;  D = "0.3130";       E = "  +0.52%  " },
    @{ Row = 10; D = "25.25";        E = "  +0.67%  " },
    @{ Row = 11; D = "0.08361";      E = "  -1.34%  " },
    @{ Row = 12; D = "1.866.69";     E = "  +1.26%  " },
    @{ Row = 13; D = "5.252";        E = "  +0.84%  " },
    @{ Row = 14; D = "0.7193";       E = $null },
    @{ Row = 15; D = "91.51";        E = "  +0.24%  " },
    @{ Row = 16; D = "6.250";        E = "  +4.36%  " },
    @{ Row = 17; D = "0.000008435";  E = "  +0.61%  " },
    @{ Row = 18; D = "29.348.86";    E = "  +0.13%  " },
    @{ Row = 19; D = "240.93";       E = "  -0.81%  " },
    @{ Row = 20; D = $null;          E = "  +0.26%  " },
    @{ Row = 21; D = "2.121.44";     E = "  -0.13%  " },
    @{ Row = 22; D = "0.9990";       E = "  -0.04%  " },
    @{ Row = 23; D = "7.806";        E = "  +0.04%  " },
    @{ Row = 24; D = "1.000";        E = "  +0.01%  " },
    @{ Row = 25; D = "0.1591";       E = "  -1.73%  " },
    @{ Row = 26; D = "163.25";       E = "  +0.25%  " },
    @{ Row = 27; D = $null;          E = "  +0.57%  " },
    @{ Row = 28; D = "18.57";        E = "  +0.33%  " },
    @{ Row = 29; D = "1.506";        E = "  -0.18%  " },
    @{ Row = 30; D = "4.423";        E = "  +0.14%  " },
    @{ Row = 31; D = "4.338";        E = "  +0.21%  " },
    @{ Row = 32; D = "1.203";        E = "  -5.67%  " },
    @{ Row = 33; D = "0.05375";      E = "  +2.27%  " },
    @{ Row = 34; D = "1.951";        E = "  +1.59%  " },
    @{ Row = 35; D = "0.7521";       E = "  +1.22%  " },
    @{ Row = 36; D = "1.179";        E = "  +0.50%  " },
    @{ Row = 37; D = "2.701";        E = "  +0.71%  " },
    @{ Row = 38; D = $null;          E = "  +1.12%  " },
    @{ Row = 39; D = "1.279.56";     E = "  +9.73%  " },
    @{ Row = 40; D = "2.739";        E = "  +0.77%  " },
    @{ Row = 41; D = "6.592";        E = "  +3.79%  " },
    @{ Row = 42; D = $null;          E = "  +3.52%  " },
    @{ Row = 43; D = "0.8916";       E = "  +0.27%  " },
    @{ Row = 44; D = "73.13";        E = $null },
    @{ Row = 45; D = $null;          E = "  +8.54%  " },
    @{ Row = 46; D = "1.000";        E = "  +0.09%  " },
    @{ Row = 47; D = "2.021.70";     E = "  +0.02%  " },
    @{ Row = 48; D = $null;          E = "  -0.36%  " },
    @{ Row = 49; D = "0.5202";       E = "  +0.17%  " },
    @{ Row = 50; D = "9.471";        E = "  +1.00%  " },
    @{ Row = 51; D = "0.4366";       E = "  +1.54%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($row, 4)
        # Force text storage so numeric-looking prices (e.g. "1.0000",
        # "0.7146") keep their exact original digits/trailing zeros
        # instead of being coerced into floating point numbers.
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
    }

    if ($null -ne $u.E) {
        # Volume strings always include surrounding spaces and a "%"
        # sign, so Excel stores them as text without any extra coercion.
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
